# "Update to add region"
#
# Adds a "Region" column to the Track BOM table on Sheet1 and introduces two
# new Imperial track rows (4ft / 8ft), in addition to the existing Metric
# (1m / 2m) rows. Also leaves Sheet1 as the active/selected sheet (it was
# MountingHardware before).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)          # Sheet1
$ws2 = $wb.Worksheets.Item(2)          # MountingHardware

# --- 1. Insert a new "Region" column at C (old C..F shift right to D..G) ---
$ws1.Columns.Item(3).Insert()
$ws1.Columns.Item(3).ColumnWidth = $ws1.Columns.Item(2).ColumnWidth

# --- 2. Insert two new rows at 4 & 5 for the Imperial track sizes ---
#     (old row 4 "End Cap" etc. shift down to row 6 onward)
$ws1.Rows.Item(4).Resize(2).Insert()

# --- 3. Fill in the Region column for the existing Track rows ---
$ws1.Range("C1").Value = "Region"
$ws1.Range("C2").Value = "Metric"
$ws1.Range("C3").Value = "Metric"

# --- 4. New row 4: Imperial 4ft track ---
$ws1.Range("A4").Value = "Imperial 4ft"
$ws1.Range("B4").Value = "Track"
$ws1.Range("C4").Value = "Imperial"
$ws1.Range("D4").Value = "ST1-9500-4FT-**"
$ws1.Range("E4").Value = "ST8-9500-4FT-**"
$ws1.Range("F4").Value = "ST3-9500-4FT-**"
$ws1.Range("G4").Value = "ST6-9500-4FT-**"

# --- 5. New row 5: Imperial 8ft track ---
$ws1.Range("A5").Value = "Imperial 8ft"
$ws1.Range("B5").Value = "Track"
$ws1.Range("C5").Value = "Imperial"
$ws1.Range("D5").Value = "ST1-9500-8FT-**"
$ws1.Range("E5").Value = "ST8-9500-8FT-**"
$ws1.Range("F5").Value = "ST3-9500-8FT-**"
$ws1.Range("G5").Value = "ST6-9500-8FT-**"

# --- 6. Make Sheet1 the active tab with G10 selected (was MountingHardware/D4) ---
$ws1.Activate() | Out-Null
$ws1.Range("G10").Select() | Out-Null
